$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("G8").Value = 3.1
$ws.Range("H8").Value = 2.32

# Row 9
$ws.Range("I9").Value = 4.3
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 3.75
$ws.Range("P9").Value = 2.34

# Row 10
$ws.Range("O10").Value = 1.32
$ws.Range("Q10").Value = 1.96
$ws.Range("AJ10").Value = 500
$ws.Range("AN10").Value = 85

# Row 12
$ws.Range("J12").Value = 3.1

# Row 13
$ws.Range("F13").Value = 2.12
$ws.Range("G13").Value = 2.3
$ws.Range("I13").Value = 4.7

# Row 14
$ws.Range("F14").Value = 2.08
$ws.Range("G14").Value = 2.2
$ws.Range("H14").Value = 4.1
$ws.Range("J14").Value = 3.1
$ws.Range("K14").Value = 3.45

# Row 16
$ws.Range("Q16").Value = 1.97
